# ECHO_CST_Xwalk.xlsx edit
#
# Summary of changes (per commit "mods to Report, server, and ui, and
# ECHO_CST_Xwalk file"):
#   1. Add a new lookup row for "Chromium"/"chromium" (STD_POLL_ID_CST 180,
#      POLLUTANT_CODE 1012), inserted right before the existing "Chrysene"
#      row (old row 64), pushing all following rows down by one.
#   2. Add a NOTES column (column E) header in E1.
#   3. Append a brand-new row at the bottom of the table for
#      "Endosulfan, total"/"endosulfan" (POLLUTANT_CODE 10757,
#      STD_POLL_ID_CST 807) with a NOTES value of "CHECK".
#   4. Keep the table's AutoFilter / _FilterDatabase defined name in sync
#      with the new data extent (header + 150 data rows = A1:E151), while
#      the sheet's overall used range grows to A1:E152 once the appended
#      row is written.
#   5. Move the active selection to G159 (matches the saved view state).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new "Chromium" row before current row 64 -----------
$ws.Rows.Item(64).Insert()

$ws.Range("A64").Value = 1012
$ws.Range("B64").Value = "Chromium"
$ws.Range("C64").Value = 180
$ws.Range("D64").Value = "chromium"

# --- 2. New NOTES header -------------------------------------------------
$ws.Range("E1").Value = "NOTES"

# --- 3. Refresh the AutoFilter / FilterDatabase to cover the header plus
#        the 150 pre-existing data rows (A1:E151) *before* appending the
#        brand-new trailing row, so the filter range doesn't balloon to
#        include it. ---------------------------------------------------
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:E151").AutoFilter()

foreach ($n in $wb.Names) {
    $n.RefersTo = "=ECHO_CST_Xwalk!`$A`$1:`$E`$151"
}

# --- 4. Append the new "Endosulfan, total" row at the end of the table --
$ws.Range("A152").Value = 10757
$ws.Range("B152").Value = "Endosulfan, total"
$ws.Range("C152").Value = 807
$ws.Range("D152").Value = "endosulfan"
$ws.Range("E152").Value = "CHECK"

# --- 5. Restore the saved selection/scroll state -------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 131
$null = $ws.Range("G159").Select()
